# Weekly update: insert a new price record as row 72 (Vega Monumental Concepción - Berenjena),
# shifting existing rows 72-116 down to 73-117.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72; this pushes old rows 72..116 down to 73..117
# and preserves formatting of the surrounding rows.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new week's data.
$ws.Cells.Item(72, 1).Value  = 11
$ws.Cells.Item(72, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(72, 3).Value  = "Bíobío"
$ws.Cells.Item(72, 4).Value  = 44813
$ws.Cells.Item(72, 5).Value  = 8
$ws.Cells.Item(72, 6).Value  = 100112001
$ws.Cells.Item(72, 7).Value  = "Berenjena"
$ws.Cells.Item(72, 8).Value  = "Sin especificar"
$ws.Cells.Item(72, 9).Value  = "Primera"
$ws.Cells.Item(72, 10).Value = 100
$ws.Cells.Item(72, 11).Value = 11000
$ws.Cells.Item(72, 12).Value = 12000
$ws.Cells.Item(72, 13).Value = 11500
$ws.Cells.Item(72, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(72, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 16).Value = 192
$ws.Cells.Item(72, 17).Value = 60
$ws.Cells.Item(72, 18).Value = "Hortaliza"
